$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "ASMADHUKUMAR"
$ws.Range("D2").Value = "CHANGESUPERVISOR"
$ws.Range("E2").Value = "PENDING"
$ws.Range("F2").Value = 1
$ws.Range("H2").Value = "YPKE"
$ws.Range("C2").Value = "BOAN"

$ws.Range("C2").Select()
